$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 223.61539
$ws.Range("I39").Value = 35.555557
$ws.Range("K39").Value = 106.666671
$ws.Range("M39").Value = 189.333329

$ws.Range("H40").Value = 7499
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7499
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7499
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7849

$ws.Range("H43").Value = 7472.0835
$ws.Range("J43").Value = 7969.5454
$ws.Range("L43").Value = 7969.5454
$ws.Range("N43").Value = -8107.5454

$ws.Range("H58").Value = 19231768
$ws.Range("I58").Value = 22727818
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 68183454
$ws.Range("L58").Value = 10500
$ws.Range("M58").Value = -68183304
$ws.Range("N58").Value = -10800

$ws.Range("H64").Value = 8207.615
$ws.Range("J64").Value = 9469.9
$ws.Range("L64").Value = 9469.9
$ws.Range("N64").Value = -9965.9

$ws.Range("H67").Value = 8207.615
$ws.Range("J67").Value = 9469.9
$ws.Range("L67").Value = 9469.9
$ws.Range("N67").Value = -11185.9

$ws.Range("H69").Value = 166677230
$ws.Range("I69").Value = 4754
$ws.Range("J69").Value = 250013470
$ws.Range("K69").Value = 14262
$ws.Range("L69").Value = 750040410
$ws.Range("M69").Value = -13388
$ws.Range("N69").Value = -750042158

$ws.Range("H72").Value = 166677230
$ws.Range("I72").Value = 4754
$ws.Range("J72").Value = 250013470
$ws.Range("K72").Value = 42786
$ws.Range("L72").Value = 2250121230
$ws.Range("M72").Value = -38418
$ws.Range("N72").Value = -2250129966

$ws.Range("H76").Value = 3105.2632

$ws.Range("H79").Value = 3105.2632

$ws.Range("H80").Value = 704.8
$ws.Range("I80").Value = 558.8571
$ws.Range("J80").Value = 783.38464
$ws.Range("K80").Value = 1676.5713
$ws.Range("L80").Value = 2350.15392
$ws.Range("M80").Value = -678.5712999999998
$ws.Range("N80").Value = -4346.15392

$ws.Range("H83").Value = 704.8
$ws.Range("I83").Value = 558.8571
$ws.Range("J83").Value = 783.38464
$ws.Range("K83").Value = 5029.7139
$ws.Range("L83").Value = 7050.46176
$ws.Range("M83").Value = -37.71389999999974
$ws.Range("N83").Value = -17034.46176

$ws.Range("H106").Value = 4252.773
$ws.Range("I106").Value = 3753.389
$ws.Range("K106").Value = 3753.389
$ws.Range("M106").Value = -3122.389

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3543.7026
$ws.Range("I32").Value = 3503.25
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 3503.25
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -3216.25
$ws.Range("N32").Value = -5574

$ws.Range("H61").Value = 4041.3809
$ws.Range("I61").Value = 4149.579
$ws.Range("K61").Value = 4149.579
$ws.Range("M61").Value = -3937.579

$ws.Range("H136").Value = 4041.3809
$ws.Range("I136").Value = 4149.579
$ws.Range("K136").Value = 12448.737
$ws.Range("M136").Value = -9898.737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3418.6155
$ws.Range("I99").Value = 2407
$ws.Range("K99").Value = 2407
$ws.Range("M99").Value = -909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1879.5555
$ws.Range("I58").Value = 1913.6875
$ws.Range("J58").Value = 1606.5
$ws.Range("K58").Value = 1913.6875
$ws.Range("L58").Value = 1606.5
$ws.Range("M58").Value = -1710.6875
$ws.Range("N58").Value = -2012.5

$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 20000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -21248

$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 100000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -106240

$ws.Range("H94").Value = 3113
$ws.Range("J94").Value = 3903.6
$ws.Range("L94").Value = 3903.6
$ws.Range("N94").Value = -4805.6

$ws.Range("H105").Value = 574.5
$ws.Range("I105").Value = 574.5
$ws.Range("K105").Value = 574.5
$ws.Range("M105").Value = 1172.5

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H136").Value = 1879.5555
$ws.Range("I136").Value = 1913.6875
$ws.Range("J136").Value = 1606.5
$ws.Range("K136").Value = 5741.0625
$ws.Range("L136").Value = 4819.5
$ws.Range("M136").Value = -3191.0625
$ws.Range("N136").Value = -9919.5

$ws.Range("H141").Value = 117995.89
$ws.Range("J141").Value = 117995.89
$ws.Range("L141").Value = 117995.89
$ws.Range("N141").Value = -128355.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.04348
$ws.Range("I2").Value = 119.916664
$ws.Range("K2").Value = 719.499984
$ws.Range("M2").Value = -606.499984

$ws.Range("H14").Value = 247.2
$ws.Range("I14").Value = 247.2
$ws.Range("K14").Value = 741.5999999999999
$ws.Range("M14").Value = -568.5999999999999

$ws.Range("H103").Value = 1499
$ws.Range("I103").Value = 1499
$ws.Range("K103").Value = 4497
$ws.Range("M103").Value = -3618

$ws.Range("H105").Value = 68000
$ws.Range("J105").Value = 68000
$ws.Range("L105").Value = 204000
$ws.Range("N105").Value = -209242

$ws.Range("H131").Value = 13335427
$ws.Range("J131").Value = 22224222
$ws.Range("L131").Value = 66672666
$ws.Range("N131").Value = -66682746

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 333333340
$ws.Range("I70").Value = 333333340
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 333333340
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -333333070
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 333333340
$ws.Range("I73").Value = 333333340
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 333333340
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -333332404
$ws.Range("N73").ClearContents()

$ws.Range("H122").Value = 17688.889
$ws.Range("I122").Value = 18560.268
$ws.Range("J122").Value = 13332
$ws.Range("K122").Value = 55680.804
$ws.Range("L122").Value = 39996
$ws.Range("M122").Value = -53230.804
$ws.Range("N122").Value = -44896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5097.92
$ws.Range("J40").Value = 7519.9
$ws.Range("L40").Value = 7519.9
$ws.Range("N40").Value = -7791.9

$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 6000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5356

$ws.Range("H122").Value = 5304.4116
$ws.Range("I122").Value = 5321.231
$ws.Range("J122").Value = 5249.75
$ws.Range("K122").Value = 15963.693
$ws.Range("L122").Value = 15749.25
$ws.Range("M122").Value = -13513.693
$ws.Range("N122").Value = -20649.25

$ws.Range("H132").Value = 3064.4688
$ws.Range("I132").Value = 2996.3914
$ws.Range("J132").Value = 3238.4443
$ws.Range("K132").Value = 8989.1742
$ws.Range("L132").Value = 9715.332900000001
$ws.Range("M132").Value = -6459.174199999999
$ws.Range("N132").Value = -14775.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 10314.833
$ws.Range("J74").Value = 10314.833
$ws.Range("L74").Value = 10314.833
$ws.Range("N74").Value = -12186.833

$ws.Range("H77").Value = 10314.833
$ws.Range("J77").Value = 10314.833
$ws.Range("L77").Value = 30944.499
$ws.Range("N77").Value = -40304.499

$ws.Range("H81").Value = 6000
$ws.Range("J81").Value = 7250
$ws.Range("L81").Value = 14500
$ws.Range("N81").Value = -16622

$ws.Range("H84").Value = 6000
$ws.Range("J84").Value = 7250
$ws.Range("L84").Value = 72500
$ws.Range("N84").Value = -83108

$ws.Range("H107").Value = 2458.1333
$ws.Range("I107").Value = 2747.6667
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 8243.000100000001
$ws.Range("L107").Value = 3900
$ws.Range("M107").Value = -6323.000100000001
$ws.Range("N107").Value = -7740

$ws.Range("H126").Value = 1890.8572
$ws.Range("I126").Value = 1813.5714
$ws.Range("J126").Value = 1968.1428
$ws.Range("K126").Value = 5440.7142
$ws.Range("L126").Value = 5904.428400000001
$ws.Range("M126").Value = -2970.7142
$ws.Range("N126").Value = -10844.4284

$ws.Range("H136").Value = 2599.4
$ws.Range("I136").Value = 1998
$ws.Range("J136").Value = 2749.75
$ws.Range("K136").Value = 5994
$ws.Range("L136").Value = 8249.25
$ws.Range("M136").Value = -3444
$ws.Range("N136").Value = -13349.25
